$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-of dates) ---
$ws.Range("A8").Value = "Volume 30   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/26/2023  Through  7/2/2023"

# --- Cells changing from blank-placeholder text to numeric (copy donor style, then set value) ---
$ws.Range("D15").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 1
$ws.Range("D15").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 2
$ws.Range("D15").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("E15").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("D15").Copy($ws.Range("C26"))
$ws.Range("C26").Value = 1
$ws.Range("D15").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$ws.Range("E15").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100

# --- Cell changing from numeric back to blank-placeholder text ---
$ws.Range("D30").Copy($ws.Range("C30"))

# --- Remaining numeric value updates ---
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 5
$ws.Range("J15").Value = 6
$ws.Range("K15").Value = -16.666666666666
$ws.Range("L15").Value = 150
$ws.Range("M15").Value = 66.666666666666
$ws.Range("N15").Value = -28.571428571428
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 75
$ws.Range("I16").Value = 47
$ws.Range("J16").Value = 41
$ws.Range("K16").Value = 14.634146341463
$ws.Range("L16").Value = 62.068965517241
$ws.Range("M16").Value = 51.612903225806
$ws.Range("N16").Value = -84.740259740259
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = -36.363636363636
$ws.Range("I17").Value = 50
$ws.Range("J17").Value = 56
$ws.Range("K17").Value = -10.714285714285
$ws.Range("L17").Value = -16.666666666666
$ws.Range("M17").Value = 47.058823529411
$ws.Range("N17").Value = -45.054945054945
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -60
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -80.95238095238
$ws.Range("I18").Value = 59
$ws.Range("J18").Value = 105
$ws.Range("K18").Value = -43.809523809523
$ws.Range("L18").Value = -3.27868852459
$ws.Range("M18").Value = 3.508771929824
$ws.Range("N18").Value = -89.350180505415
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -30.76923076923
$ws.Range("F19").Value = 56
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = 9.803921568627
$ws.Range("I19").Value = 368
$ws.Range("J19").Value = 320
$ws.Range("K19").Value = 15
$ws.Range("L19").Value = 53.974895397489
$ws.Range("M19").Value = 4.545454545454
$ws.Range("N19").Value = -66.207529843893
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -63.636363636363
$ws.Range("I20").Value = 26
$ws.Range("J20").Value = 27
$ws.Range("K20").Value = -3.703703703703
$ws.Range("L20").Value = -13.333333333333
$ws.Range("M20").Value = 85.714285714285
$ws.Range("N20").Value = -91.849529780564
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -45.16129032258
$ws.Range("F21").Value = 88
$ws.Range("G21").Value = 104
$ws.Range("H21").Value = -15.384615384615
$ws.Range("I21").Value = 555
$ws.Range("J21").Value = 555
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 31.828978622327
$ws.Range("M21").Value = 13.034623217922
$ws.Range("N21").Value = -76.5625
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("J22").Value = 11
$ws.Range("K22").Value = 45.454545454545
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = -25.714285714285
$ws.Range("F24").Value = 90
$ws.Range("G24").Value = 112
$ws.Range("H24").Value = -19.642857142857
$ws.Range("I24").Value = 489
$ws.Range("J24").Value = 623
$ws.Range("K24").Value = -21.508828250401
$ws.Range("L24").Value = -2.2
$ws.Range("M24").Value = 53.291536050156
$ws.Range("C25").Value = 4
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 16
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 134
$ws.Range("J25").Value = 118
$ws.Range("K25").Value = 13.559322033898
$ws.Range("L25").Value = 78.666666666666
$ws.Range("M25").Value = 8.064516129032
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 7
$ws.Range("J26").Value = 8
$ws.Range("K26").Value = -12.5
$ws.Range("L26").Value = 40
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 33.333333333333
$ws.Range("J27").Value = 24
$ws.Range("K27").Value = 25
$ws.Range("L27").Value = -14.285714285714
